# Update "Des Scheduled Flights vs actual.xlsx":
#  - append 17 new daily rows (2021-03-01 .. 2021-03-17) to the Ark1 sheet
#  - extend the shared % formula down to the new rows
#  - move the active selection / view roughly to where the author left it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to append: Date (text), Scheduled flights, Tracked flights
$newData = @(
    ,@("2021-03-01", 52, 50)
    ,@("2021-03-02", 53, 52)
    ,@("2021-03-03", 50, 47)
    ,@("2021-03-04", 58, 57)
    ,@("2021-03-05", 68, 64)
    ,@("2021-03-06", 42, 41)
    ,@("2021-03-07", 47, 47)
    ,@("2021-03-08", 56, 55)
    ,@("2021-03-09", 49, 44)
    ,@("2021-03-10", 57, 56)
    ,@("2021-03-11", 60, 59)
    ,@("2021-03-12", 81, 77)
    ,@("2021-03-13", 56, 52)
    ,@("2021-03-14", 53, 45)
    ,@("2021-03-15", 47, 39)
    ,@("2021-03-16", 49, 46)
    ,@("2021-03-17", 59, 58)
)

$firstNewRow = 329
$lastExisting = 328
$lastNewRow = $firstNewRow + $newData.Count - 1

# Carry the existing formatting (text date column, number columns, percent
# column) down onto the freshly-appended rows before writing any values.
$srcFormatRow = $ws.Range("A$lastExisting`:D$lastExisting")
$dstFormatRange = $ws.Range("A$firstNewRow`:D$lastNewRow")
$srcFormatRow.Copy()
$dstFormatRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write the new rows.
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $firstNewRow + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Fill the percent formula down through the new rows.
$ws.Range("D$firstNewRow`:D$lastNewRow").Formula = "=C$firstNewRow/B$firstNewRow"

# Recalculate so the cached <v> values land correctly.
$excel.Calculate()

# Move the view roughly back to where the author left off and restore the
# window geometry recorded in the workbook.
$win = $excel.ActiveWindow
$win.ScrollRow = 187
$win.ScrollColumn = 1
$ws.Range("O203").Select()
$win.Left = 2580
$win.Top = 3240
